## Adds a new "wbCloseAppsRecover" workblock (RecoverApps) to the Workblocks sheet,
## right after the Init workblock, and shortens the existing workblock "type" values
## (which used to hold long "GetData, Framework/Business, X" strings) down to their
## short form, matching the updated wbLogging interface (wbKey / wbPath).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workblocks")

# Insert two new rows for the new "RecoverApps" workblock, right after the
# existing "Init" workblock entry (rows 3-4).
$ws.Range("A5:A6").EntireRow.Insert()

$ws.Range("A5").Value = "wbCloseAppsRecover_Type"
$ws.Range("C5").Value = "Name of Workblock"

$ws.Range("A6").Value = "wbCloseAppsRecover_SuppressSuccessful"
$ws.Range("B6").Value = $true
$ws.Range("C6").Value = "Do not log successful executions of wb"

# Shorten all the workblock "type" values (previously long "GetData, X, Y"
# style strings) down to just the short name used by the new interface.
$ws.Range("B3").Value = "Init"
$ws.Range("B5").Value = "RecoverApps"
$ws.Range("B7").Value = "GetData"
$ws.Range("B9").Value = "Process"
$ws.Range("B11").Value = "Next"
$ws.Range("B13").Value = "CloseApps"
$ws.Range("B15").Value = "InitApps"
$ws.Range("B17").Value = "ProcessApps"

# Match the workbook's row formatting exactly.
$ws.Range("B13").Style = "Normal"
$ws.Range("B15").HorizontalAlignment = -4131
$ws.Range("B17").Style = "Normal"
$ws.Range("B18").Style = "Normal"

# Re-select the full table, matching the updated selection in the workbook.
$ws.Activate()
$ws.Range("A3:C18").Select()
